# Deploy the implementation guide.
#
# Updates the "Metadata" sheet of the CodeSystem-family-type workbook:
#   - Status  (B6): active -> draft
#   - Date    (B8): 2023-05-12T12:33:13+00:00 -> 2023-08-01T16:12:28+00:00
# and (re)applies the vertical-top / wrap-text alignment that the header
# row and data rows already use. The header style and the body style are
# shared across both worksheets, so the alignment is reasserted on both
# the "Metadata" and "Concepts" sheets.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# -- Status: active -> draft --------------------------------------------
$ws.Range("B6").Value = "draft"

# -- Date: regenerated on deploy -----------------------------------------
$ws.Range("B8").Value = "2023-08-01T16:12:28+00:00"

# -- Re-affirm the workbook's alignment (vertical top + wrap text) so the
#    header row and the data rows keep their intended formatting explicit
#    on every sheet that shares those cell styles.
foreach ($sheetName in @("Metadata", "Concepts")) {
    $sheet = $wb.Worksheets.Item($sheetName)
    $used = $sheet.UsedRange

    $headerRow = $used.Rows.Item(1)
    $headerRow.VerticalAlignment = -4160   # xlTop
    $headerRow.WrapText = $true

    if ($used.Rows.Count -gt 1) {
        $bodyRows = $used.Rows.Item(2).Resize($used.Rows.Count - 1)
        $bodyRows.VerticalAlignment = -4160   # xlTop
        $bodyRows.WrapText = $true
    }
}
